$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 11
$ws.Range("Q5").Value = 1.88
$ws.Range("R5").Value = 1.98
$ws.Range("P6").Value = 4.05
$ws.Range("U6").Value = 1.85
$ws.Range("V6").Value = 1.91
$ws.Range("J7").Value = 3.15
$ws.Range("K7").Value = 2.07
$ws.Range("L7").Value = 3.15
$ws.Range("O7").Value = 1.24
$ws.Range("P7").Value = 3.3
$ws.Range("Q7").Value = 1.7
$ws.Range("R7").Value = 1.91
$ws.Range("S7").Value = 1.34
$ws.Range("T7").Value = 3.1
$ws.Range("U7").Value = 1.53
$ws.Range("V7").Value = 2.18
$ws.Range("W7").Value = 9.75
$ws.Range("X7").Value = 14.5
$ws.Range("AB7").Value = 25
$ws.Range("AC7").Value = 10.75
$ws.Range("AE7").Value = 11.25
$ws.Range("AH7").Value = 10
$ws.Range("AI7").Value = 15
$ws.Range("AL7").Value = 20
$ws.Range("AM7").Value = 24
$ws.Range("AO7").Value = 14
$ws.Range("AP7").Value = 19
$ws.Range("AR7").Value = 80
$ws.Range("AT7").Value = 2.82
$ws.Range("AX7").Value = 14
$ws.Range("AY7").Value = 18.5
$ws.Range("BB7").Value = 200
$ws.Range("G11").Value = 2
$ws.Range("I11").Value = 3.8
$ws.Range("N11").Value = 7.5
$ws.Range("O11").Value = 1.44
$ws.Range("P11").Value = 2.63
$ws.Range("Q11").Value = 2.35
$ws.Range("R11").Value = 1.57
$ws.Range("AH11").Value = 9
$ws.Range("AN11").Value = 3.75
$ws.Range("AO11").Value = 11
$ws.Range("N12").Value = 8
